# Updates the A02 Savings Account unit-test plan:
#  - Rewrites the Preconditions / Method Inputs / Expected Result columns
#    (E:G) for test cases in rows 7-12 with the revised wording from the
#    finished test plan.
#  - Moves the active selection from G12 to E12 (matches the saved view
#    state captured in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: __init__ / valid data
$ws.Range("E7").Value = "Valid account data"
$ws.Range("F7").Value = '("5001", "1003", 150, date(2022, 5, 10), 50.0)'
$ws.Range("G7").Value = 'account_number="5001", balance=150, date_created=date(2022, 5, 10), minimum_balance=50.0'

# Row 8: __init__ / invalid minimum_balance
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = '("5002", "1004", 200, date(2023, 6, 10), "invalid")'
$ws.Range("G8").Value = "minimum_balance defaults to 50.0"

# Row 9: get_service_charges / balance greater than minimum balance
$ws.Range("E9").Value = "self.balance = 200"
$ws.Range("F9").Value = "get_service_charges()"
$ws.Range("G9").Value = "Returns base service charge of 0.50"

# Row 10: get_service_charges / balance equal to minimum balance
$ws.Range("E10").Value = "self.balance = 50"
$ws.Range("F10").Value = "get_service_charges()"
$ws.Range("G10").Value = "Returns base service charge of 0.50"

# Row 11: get_service_charges / balance less than minimum balance
$ws.Range("E11").Value = "self.balance = 20"
$ws.Range("F11").Value = "get_service_charges()"
$ws.Range("G11").Value = "Returns 1.00 (0.50 * 2.0 due to SERVICE_CHARGE_PREMIUM)"

# Row 12: __str__
$ws.Range("E12").Value = "Instance initialized with valid attributes"
$ws.Range("F12").Value = "str(savings_account)"
$ws.Range("G12").Value = "Returns ""Account Number: 5001 Balance: `$150.00\nMinimum Balance: `$50.00 Account Type: Savings"""

# Move the saved selection from G12 to E12.
$ws.Range("E12").Select()
